# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (with fund holding data) right after the
# "总计" (totals) sheet, and adds a corresponding summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as TEXT (inline/shared string), even when
# the value looks numeric (e.g. "3.22", "007731"), without leaving a stray
# NumberFormat style behind on the cell.
# ---------------------------------------------------------------------------
function Set-TextCell {
    param($Cell, $Val)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Val
    $Cell.ClearFormats()
}

# NOTE: this runtime's PowerShell engine does not reliably bind *named*
# parameters (e.g. `-Cell $c -Val $v` silently yields $null inside the
# function). Always call Set-TextCell positionally: Set-TextCell $cell $val

# ===========================================================================
# 1. Insert the new "2022-Q3" worksheet right after "总计"
# ===========================================================================
$totalSheet = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Header row (B1:H1)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}
# Match the bold/centered/bordered header style used by the other sheets
$totalSheet.Cells.Item(1, 2).Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

# Data rows
$q3Data = @(
    @("007731", "民生加银持续成长混合A",   "3.22", "94.57", "5.69", "0.1832", 7),
    @("007732", "民生加银持续成长混合C",   "1.89", "94.57", "5.69", "0.1075", 7),
    @("163818", "中银中小盘成长混合",       "0.77", "88.59", "2.60", "0.0200", 7),
    @("002631", "江信瑞福灵活配置混合C",   "0.39", "84.78", "4.81", "0.0188", 6),
    @("519097", "新华中小市值优选混合",     "0.71", "67.35", "2.58", "0.0183", 10),
    @("002630", "江信瑞福灵活配置混合A",   "0.01", "84.78", "4.81", "0.0005", 6)
)

for ($r = 0; $r -lt $q3Data.Length; $r++) {
    $row = $r + 2
    $rec = $q3Data[$r]

    $aCell = $q3.Cells.Item($row, 1)
    $aCell.Value = $r
    Set-TextCell $q3.Cells.Item($row, 2) $rec[0]
    Set-TextCell $q3.Cells.Item($row, 3) $rec[1]
    Set-TextCell $q3.Cells.Item($row, 4) $rec[2]
    Set-TextCell $q3.Cells.Item($row, 5) $rec[3]
    Set-TextCell $q3.Cells.Item($row, 6) $rec[4]
    Set-TextCell $q3.Cells.Item($row, 7) $rec[5]
    $q3.Cells.Item($row, 8).Value = $rec[6]
}

# Match the column-A index style (bold/centered/bordered) used elsewhere
$totalSheet.Cells.Item(1, 2).Copy()
$q3.Range("A2:A7").PasteSpecial(-4122)

# ===========================================================================
# 2. Add the "2022-Q3" summary row to "总计" (shifting Q2/Q4 rows down)
# ===========================================================================
$b2 = $totalSheet.Cells.Item(2, 2).Value()
$c2 = $totalSheet.Cells.Item(2, 3).Value()
$d2 = $totalSheet.Cells.Item(2, 4).Value()
$b3 = $totalSheet.Cells.Item(3, 2).Value()
$c3 = $totalSheet.Cells.Item(3, 3).Value()
$d3 = $totalSheet.Cells.Item(3, 4).Value()

# old row3 (2021-Q4) -> row4
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = $b3
$totalSheet.Cells.Item(4, 3).Value = $c3
$totalSheet.Cells.Item(4, 4).Value = $d3

# old row2 (2022-Q2) -> row3
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = $b2
$totalSheet.Cells.Item(3, 3).Value = $c2
$totalSheet.Cells.Item(3, 4).Value = $d2

# new row2 (2022-Q3)
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 0.35

# Fix up the column-A style on the newly created row 4 (copy from row 3)
$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(4, 1).PasteSpecial(-4122)

Write-Host "2022-Q3 sheet added; 总计 updated"
